# feat: add 2022-Q3 data
#
# Adds a new "2022-Q3" sheet (same fund-holdings layout as the other
# quarters, cloned from "2021-Q4" for formatting) and inserts the
# matching summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

$totalSheet    = $wb.Worksheets.Item("总计")
$q2Sheet       = $wb.Worksheets.Item("2022-Q2")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# Helper: write $text into $destRange as literal text, even when it
# looks numeric (e.g. "519677" or "2.88"), without leaving a custom
# number-format/style on the destination cell. Excel normally "smart"
# converts numeric-looking strings assigned via .Value into numbers, so
# we stage the text (quote-prefixed, forcing text) in a scratch cell and
# PasteSpecial only the *value* across - that carries the text/shared-
# string type over but none of the scratch cell's formatting.
function Set-TextValue($ws, $destRange, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $destRange.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# --- 1. Create the new "2022-Q3" sheet right before "2022-Q2" ----------
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Clone the single-fund-row layout (headers + formatting) from 2021-Q4,
# then overwrite with the new quarter's values.
$templateSheet.Range("B1:H2").Copy($q3Sheet.Range("B1"))
$templateSheet.Range("A2").Copy($q3Sheet.Range("A2"))

Set-TextValue $q3Sheet $q3Sheet.Range("B2") "519677"
$q3Sheet.Range("C2").Value = "银河定投宝腾讯济安指数"
Set-TextValue $q3Sheet $q3Sheet.Range("D2") "2.88"
Set-TextValue $q3Sheet $q3Sheet.Range("E2") "92.40"
Set-TextValue $q3Sheet $q3Sheet.Range("F2") "1.26"
Set-TextValue $q3Sheet $q3Sheet.Range("G2") "0.0363"
$q3Sheet.Range("H2").Value = 8

# --- 2. Insert the new summary row on "总计" above the old 2022-Q2 row -
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.04

# Renumber the index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
